$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - D2 description text updated
$ws.Range("D2").Value = 'Это приложение к диплому, где указывается специальность и степень, которую получает студент (в данном случае - бакалавр-инженер в области Информатики и Вычислительная Техника).'

# Row 3 - C3, D3, E3, F3 updated
$ws.Range("C3").Value = 'Диплом бакалавра'
$ws.Range("D3").Value = 'Этот документ присвоен степени бакалавра и содержит сведения о личности обладателя диплома, ученом звании, специальности, дате рождения, предыдущем уровне образования и названию и адресу учебного заведения.'
$ws.Range("E3").Value = '✅'
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '1.00'
$ws.Range("F3").Style = "Normal"

# Row 4 - C4, D4, E4, F4 updated
$ws.Range("C4").Value = 'Приложение к диплому'
$ws.Range("D4").Value = 'Содержит список дисциплин (модулей) основной образовательной программы в высшем профессиональном учебном заведении.'
$ws.Range("E4").Value = '✅'
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '1.00'
$ws.Range("F4").Style = "Normal"

# Row 5 - D5 updated
$ws.Range("D5").Value = 'Этот документ подтверждает результаты теста по английскому языку (IELTS), который был проведён для определения способности кандидата к обучению на английском языке.'
